# Update "paises.xlsx" (sheet "Pais") with the refreshed COVID-19 snapshot:
# the table is ranked by "Casos totales" (column B) descending, so as counts
# changed several countries shifted rank/row; re-point column A's country
# name per row accordingly, refresh the stat columns that actually received
# new numbers, and bump the "datos actualizados" timestamp in A1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update country names (column A) whose rank position changed ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 29 de Marzo de 2020 a las 06:20"
$ws.Cells.Item(135, 1).Value = "Polinesia Francesa"
$ws.Cells.Item(136, 1).Value = "Uganda"
$ws.Cells.Item(149, 1).Value = "Haiti"
$ws.Cells.Item(150, 1).Value = "Republica de Yibuti"
$ws.Cells.Item(151, 1).Value = "Tanzania"
$ws.Cells.Item(153, 1).Value = "Guinea Ecuatorial"
$ws.Cells.Item(154, 1).Value = "Dominica"
$ws.Cells.Item(155, 1).Value = "San Martin (Parte Francesa)"
$ws.Cells.Item(156, 1).Value = "Bahamas"
$ws.Cells.Item(157, 1).Value = "Niger"
$ws.Cells.Item(158, 1).Value = "Groenlandia"
$ws.Cells.Item(159, 1).Value = "Suazilandia"
$ws.Cells.Item(160, 1).Value = "Birmania"
$ws.Cells.Item(161, 1).Value = "Laos"
$ws.Cells.Item(162, 1).Value = "Seychelles"
$ws.Cells.Item(163, 1).Value = "Surinam"
$ws.Cells.Item(164, 1).Value = "Mozambique"
$ws.Cells.Item(165, 1).Value = "Guinea"
$ws.Cells.Item(170, 1).Value = "Antigua y Barbuda"
$ws.Cells.Item(171, 1).Value = "Granada"
$ws.Cells.Item(174, 1).Value = "Eritrea"
$ws.Cells.Item(175, 1).Value = "Benin"
$ws.Cells.Item(176, 1).Value = "Santa Sede"
$ws.Cells.Item(178, 1).Value = "Montserrat"
$ws.Cells.Item(179, 1).Value = "Fiyi"
$ws.Cells.Item(181, 1).Value = "Siria"
$ws.Cells.Item(182, 1).Value = "Mauritania"
$ws.Cells.Item(183, 1).Value = "San Bartolome"
$ws.Cells.Item(189, 1).Value = "Republica de Africa Central"
$ws.Cells.Item(190, 1).Value = "Butan"
$ws.Cells.Item(191, 1).Value = "Liberia"
$ws.Cells.Item(192, 1).Value = "Libia"
$ws.Cells.Item(193, 1).Value = "Republica del Chad"
$ws.Cells.Item(194, 1).Value = "San Martin (Parte Holandesa)"
$ws.Cells.Item(196, 1).Value = "Santa Lucia"
$ws.Cells.Item(197, 1).Value = "Gambia"
$ws.Cells.Item(198, 1).Value = "Anguila"
$ws.Cells.Item(200, 1).Value = "Belice"
$ws.Cells.Item(202, 1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(203, 1).Value = "Timor Oriental"
$ws.Cells.Item(204, 1).Value = "Papua Nueva Guinea"

# --- Update numeric stats for rows with new/changed data ---
$ws.Cells.Item(127, 2).Value = 43
$ws.Cells.Item(127, 8).Value = 1
$ws.Cells.Item(149, 2).Value = 15
$ws.Cells.Item(149, 3).Value = 7
$ws.Cells.Item(149, 4).Value = 1
$ws.Cells.Item(150, 4).Value = 0
$ws.Cells.Item(150, 5).Value = 14
$ws.Cells.Item(151, 2).Value = 14
$ws.Cells.Item(151, 4).Value = 1
$ws.Cells.Item(151, 5).Value = 13
$ws.Cells.Item(153, 2).Value = 12
$ws.Cells.Item(153, 5).Value = 12
$ws.Cells.Item(155, 2).Value = 11
$ws.Cells.Item(155, 4).Value = 0
$ws.Cells.Item(155, 5).Value = 11
$ws.Cells.Item(156, 4).Value = 1
$ws.Cells.Item(156, 8).Value = 0
$ws.Cells.Item(157, 4).Value = 0
$ws.Cells.Item(157, 5).Value = 9
$ws.Cells.Item(157, 8).Value = 1
$ws.Cells.Item(158, 2).Value = 10
$ws.Cells.Item(158, 4).Value = 2
$ws.Cells.Item(158, 5).Value = 8
$ws.Cells.Item(159, 2).Value = 9
$ws.Cells.Item(159, 5).Value = 9
$ws.Cells.Item(196, 4).Value = 1
$ws.Cells.Item(196, 8).Value = 0
$ws.Cells.Item(197, 4).Value = 0
$ws.Cells.Item(197, 8).Value = 1
